$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions): update "想去人数" (F column) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 13
$ws1.Range("F4").Value = 3338
$ws1.Range("F7").Value = 162

# Sheet "演出" (performances): update "想去人数" (F column) value
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 121

# Sheet "全部类型" (all types): update "想去人数" (F column) values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 121
$ws4.Range("F7").Value = 13
$ws4.Range("F8").Value = 3338
$ws4.Range("F12").Value = 162
